$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 307. This pushes the
# existing rows 307..377 down to 308..378, which is exactly what the
# target diff shows (every row's data moved down by one row).
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with its data.
$ws.Cells.Item(307, 1).Value = 5
$ws.Cells.Item(307, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(307, 3).Value = "Maule"
$ws.Cells.Item(307, 4).Value = 44511
$ws.Cells.Item(307, 5).Value = 7
$ws.Cells.Item(307, 6).Value = 100112002
$ws.Cells.Item(307, 7).Value = "Pimiento"
$ws.Cells.Item(307, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 200
$ws.Cells.Item(307, 11).Value = 25000
$ws.Cells.Item(307, 12).Value = 25000
$ws.Cells.Item(307, 13).Value = 25000
$ws.Cells.Item(307, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(307, 15).Value = "Región del Maule"
$ws.Cells.Item(307, 16).Value = 1667
$ws.Cells.Item(307, 17).Value = 15
$ws.Cells.Item(307, 18).Value = "Hortaliza"
